$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# Row 1 (email header stays the same; hyperlink text changes)
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "danielvisca96@gmail.com"

# Row 2 (shop header row is unchanged)
$ws.Range("A2").Value = "shop"
$ws.Range("B2").Value = "Jubilant Jelly"

# Row 3
$ws.Range("A3").Value = "Almond Butter"
$ws.Range("B3").Value = 100

# Row 4
$ws.Range("A4").Value = "Alan Watts"
$ws.Range("B4").Value = 2

# Row 5 (new)
$ws.Range("A5").Value = "Crepe"
$ws.Range("B5").Value = 7

# Row 6 (new content, was previously "Peanut Butter"/150) - no longer uses
# the special font style that the old A6 had, so reset it to Normal.
$ws.Range("A6").Value = "Jelly Fish Sandwich"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 99

# Row 7 stays blank (gap row)

# Row 8 (new "shop"/"Jubilant Jam" pair)
$ws.Range("A8").Value = "shop"
$ws.Range("B8").Value = "Jubilant Jam"

# Row 9 (new)
$ws.Range("A9").Value = "The Answer to Life, The Universe and Everything"
$ws.Range("B9").Value = 42

# Row 10 (new)
$ws.Range("A10").Value = "Unjammer"
$ws.Range("B10").Value = 1

# Row 11 (new; final "Peanut Butter"/3 entry)
$ws.Range("A11").Value = "Peanut Butter"
$ws.Range("B11").Value = 3

# --- Column widths -----------------------------------------------------
# (ColumnWidth is stored internally on a 1/6-character grid, so these are
# the closest inputs that round-trip to the target stored widths of
# ~39.8203125 / ~23.9375 characters.)
$ws.Columns.Item(1).ColumnWidth = 39
$ws.Columns.Item(2).ColumnWidth = 23.15

# --- View / selection -----------------------------------------------------
$excel.ActiveWindow.Zoom = 105
$null = $ws.Range("B11").Select()
